$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 275
$ws.Range("I9").Value = 275
$ws.Range("K9").Value = 275
$ws.Range("M9").Value = -106
$ws.Range("H34").Value = 16331.546
$ws.Range("I34").Value = 4405.222
$ws.Range("J34").Value = 70000
$ws.Range("K34").Value = 4405.222
$ws.Range("L34").Value = 70000
$ws.Range("M34").Value = -4202.222
$ws.Range("N34").Value = -70406
$ws.Range("H36").Value = 16331.546
$ws.Range("I36").Value = 4405.222
$ws.Range("J36").Value = 70000
$ws.Range("K36").Value = 4405.222
$ws.Range("L36").Value = 70000
$ws.Range("M36").Value = -3690.222
$ws.Range("N36").Value = -71430
$ws.Range("H51").Value = 11999.4
$ws.Range("I51").Value = 9999
$ws.Range("J51").Value = 15000
$ws.Range("K51").Value = 9999
$ws.Range("L51").Value = 15000
$ws.Range("M51").Value = -9515
$ws.Range("N51").Value = -15968
$ws.Range("H70").Value = 3226.9092
$ws.Range("I70").Value = 2749.5
$ws.Range("K70").Value = 8248.5
$ws.Range("M70").Value = -7978.5
$ws.Range("H73").Value = 3226.9092
$ws.Range("I73").Value = 2749.5
$ws.Range("K73").Value = 8248.5
$ws.Range("M73").Value = -7312.5
$ws.Range("H92").Value = 66667092
$ws.Range("I92").Value = 83333760
$ws.Range("K92").Value = 83333760
$ws.Range("M92").Value = -83332512
$ws.Range("H99").Value = 66666856
$ws.Range("I99").Value = 66666856
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 200000568
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -199999070
$ws.Range("N99").ClearContents()
$ws.Range("H106").Value = 142858750
$ws.Range("I106").Value = 200001550
$ws.Range("J106").Value = 1750
$ws.Range("K106").Value = 200001550
$ws.Range("L106").Value = 1750
$ws.Range("M106").Value = -200000919
$ws.Range("N106").Value = -3012
$ws.Range("H111").Value = 1886.4117
$ws.Range("J111").Value = 1199.8334
$ws.Range("L111").Value = 3599.5002
$ws.Range("N111").Value = -9733.5002
$ws.Range("H116").Value = 4005
$ws.Range("J116").Value = 4005
$ws.Range("L116").Value = 4005
$ws.Range("N116").Value = -10889

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 10000
$ws.Range("J6").Value = 10000
$ws.Range("L6").Value = 10000
$ws.Range("N6").Value = -10346
$ws.Range("H12").Value = 571.4286
$ws.Range("I12").Value = 500
$ws.Range("K12").Value = 500
$ws.Range("M12").Value = -327
$ws.Range("H26").Value = 3833.3333
$ws.Range("I26").Value = 3000
$ws.Range("K26").Value = 3000
$ws.Range("M26").Value = -2670
$ws.Range("H30").Value = 2505987.2
$ws.Range("I30").Value = 4287071
$ws.Range("J30").Value = 12470
$ws.Range("K30").Value = 4287071
$ws.Range("L30").Value = 12470
$ws.Range("M30").Value = -4286921
$ws.Range("N30").Value = -12770
$ws.Range("H61").Value = 3800
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H97").Value = 27779546
$ws.Range("I97").Value = 41668824
$ws.Range("J97").Value = 993.75
$ws.Range("K97").Value = 41668824
$ws.Range("L97").Value = 993.75
$ws.Range("M97").Value = -41668328
$ws.Range("N97").Value = -1985.75
$ws.Range("H132").Value = 4031
$ws.Range("I132").Value = 4031
$ws.Range("K132").Value = 12093
$ws.Range("M132").Value = -9563
$ws.Range("H136").Value = 3800
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3666.6667
$ws.Range("I20").Value = 3750
$ws.Range("K20").Value = 3750
$ws.Range("M20").Value = -3503
$ws.Range("H37").Value = 2555
$ws.Range("I37").Value = 710
$ws.Range("K37").Value = 710
$ws.Range("M37").Value = -573
$ws.Range("H54").Value = 3000
$ws.Range("I54").Value = 3000
$ws.Range("K54").Value = 3000
$ws.Range("M54").Value = -2516
$ws.Range("H75").Value = 5154.3335
$ws.Range("I75").Value = 5154.3335
$ws.Range("K75").Value = 5154.3335
$ws.Range("M75").Value = -4218.3335
$ws.Range("H78").Value = 5154.3335
$ws.Range("I78").Value = 5154.3335
$ws.Range("K78").Value = 15463.0005
$ws.Range("M78").Value = -10783.0005
$ws.Range("H105").Value = 2200
$ws.Range("I105").Value = 2200
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2200
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -453
$ws.Range("N105").ClearContents()
$ws.Range("H134").Value = 901.4
$ws.Range("I134").Value = 901.4
$ws.Range("K134").Value = 2704.2
$ws.Range("M134").Value = -169.1999999999998

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1758.1428
$ws.Range("I58").Value = 1717.8334
$ws.Range("K58").Value = 1717.8334
$ws.Range("M58").Value = -1514.8334
$ws.Range("H88").Value = 34762.43
$ws.Range("J88").Value = 34762.43
$ws.Range("L88").Value = 34762.43
$ws.Range("N88").Value = -35574.43
$ws.Range("H91").Value = 34762.43
$ws.Range("J91").Value = 34762.43
$ws.Range("L91").Value = 34762.43
$ws.Range("N91").Value = -37570.43
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H134").Value = 2963.1667
$ws.Range("I134").Value = 3259.6667
$ws.Range("K134").Value = 9779.000100000001
$ws.Range("M134").Value = -7244.000100000001
$ws.Range("H136").Value = 1758.1428
$ws.Range("I136").Value = 1717.8334
$ws.Range("K136").Value = 5153.5002
$ws.Range("M136").Value = -2603.5002

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 418.83334
$ws.Range("I108").Value = 418.83334
$ws.Range("K108").Value = 1256.50002
$ws.Range("M108").Value = 1623.49998
$ws.Range("H109").Value = 450
$ws.Range("I109").Value = 450
$ws.Range("K109").Value = 1350
$ws.Range("M109").Value = -310
$ws.Range("H121").Value = 390.91666
$ws.Range("I121").Value = 195.44444
$ws.Range("K121").Value = 586.33332
$ws.Range("M121").Value = 723.66668

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 97.8
$ws.Range("I13").Value = 97.8
$ws.Range("K13").Value = 97.8
$ws.Range("M13").Value = 41.2
$ws.Range("H45").Value = 40000
$ws.Range("I45").Value = 40000
$ws.Range("K45").Value = 40000
$ws.Range("M45").Value = -39441
$ws.Range("H51").Value = 60000
$ws.Range("J51").Value = 60000
$ws.Range("L51").Value = 60000
$ws.Range("N51").Value = -61018
$ws.Range("H102").Value = 1013.2857
$ws.Range("I102").Value = 613.2
$ws.Range("K102").Value = 613.2
$ws.Range("M102").Value = 1008.8
$ws.Range("H126").Value = 1598
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 2500
$ws.Range("I132").Value = 2500
$ws.Range("K132").Value = 7500
$ws.Range("M132").Value = -4970

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 218.41667
$ws.Range("I22").Value = 199.33333
$ws.Range("K22").Value = 199.33333
$ws.Range("M22").Value = 95.66667000000001
$ws.Range("H27").Value = 218.41667
$ws.Range("I27").Value = 199.33333
$ws.Range("K27").Value = 199.33333
$ws.Range("M27").Value = -92.33332999999999
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H55").Value = 1040.05
$ws.Range("I55").Value = 701.8333
$ws.Range("K55").Value = 701.8333
$ws.Range("M55").Value = -528.8333
$ws.Range("H76").Value = 10000
$ws.Range("I76").Value = 10000
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 10000
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("M76").Value = -9662
$ws.Range("H79").Value = 10000
$ws.Range("I79").Value = 10000
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 10000
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("M79").Value = -8830
$ws.Range("H104").Value = 43999.75
$ws.Range("J104").Value = 43999.75
$ws.Range("L104").Value = 43999.75
$ws.Range("N104").Value = -50987.75
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H69").Value = 15044.833
$ws.Range("J69").Value = 15044.833
$ws.Range("L69").Value = 15044.833
$ws.Range("N69").Value = -16542.833
$ws.Range("H72").Value = 15044.833
$ws.Range("J72").Value = 15044.833
$ws.Range("L72").Value = 45134.499
$ws.Range("N72").Value = -52622.499
$ws.Range("H113").Value = 445.82352
$ws.Range("I113").Value = 251.88889
$ws.Range("K113").Value = 755.6666700000001
$ws.Range("M113").Value = 1414.33333
$ws.Range("H122").Value = 1451
$ws.Range("I122").Value = 1001.5
$ws.Range("K122").Value = 3004.5
$ws.Range("M122").Value = -554.5
